$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "想去人数" values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 8078
$wsExhibition.Range("F13").Value = 461
$wsExhibition.Range("F16").Value = 32
$wsExhibition.Range("F17").Value = 5996
$wsExhibition.Range("F20").Value = 2107
$wsExhibition.Range("F24").Value = 426

# Sheet "全部类型" (sheet4): update F column "想去人数" values (same events, different rows)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8078
$wsAll.Range("F14").Value = 461
$wsAll.Range("F17").Value = 32
$wsAll.Range("F19").Value = 5996
$wsAll.Range("F23").Value = 2107
$wsAll.Range("F27").Value = 426
